$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Historias de Usuario")
$ws.Activate()

$ws.Range("B5").Value = "H. usuario #1, #14 y  H. técnica #1"
$ws.Range("C5").Value = "Como un usuario nuevo, necesito registrarme con mi email y contraseña, con la finalidad de acceder a la red social"
$ws.Range("D5").Value = "Registro de usuarios"
$ws.Range("E5").Value = "planificada"
$ws.Range("F5").Value = "20 horas"
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = "Alta"

$ws.Range("B6").Value = "H. usuario #2 y H. técnica #2"
$ws.Range("C6").Value = "Como un usuario registrado, necesito ver y editar mis datos de perfil, con la finalidad de tener mi información actualizada"
$ws.Range("D6").Value = "Gestión de perfil de usuarios"
$ws.Range("E6").Value = "planificada"
$ws.Range("F6").Value = "16 horas"
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = "Alta"

$ws.Range("B7").Value = "H. usuario #3 y H. técnica #3"
$ws.Range("C7").Value = "Como un usuario activo, necesito publicar contenido en el feed, con la finalidad de interactuar con otros usuarios y compartir mis pensamientos"
$ws.Range("D7").Value = "Publicaciones de contenido"
$ws.Range("E7").Value = "planificada"
$ws.Range("F7").Value = "22 horas"
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = "Alta"

$ws.Range("B8").Value = "H. usuario #4 y H. técnica #4"
$ws.Range("C8").Value = "Como un usuario, necesito generar imágenes a partir de un prompt, con la finalidad de agregar contenido visual único a mis publicaciones"
$ws.Range("D8").Value = "Generación de imágenes con IA"
$ws.Range("E8").Value = "planificada"
$ws.Range("F8").Value = "24 horas"
$ws.Range("G8").Value = 3
$ws.Range("H8").Value = "Alta"

$ws.Range("B9").Value = "H. usuario #5 y H. técnica #5"
$ws.Range("C9").Value = "Como un usuario, necesito comentar y reaccionar a publicaciones, con la finalidad de interactuar con los demás usuarios"
$ws.Range("D9").Value = "Interacciones sociales (reacciones y comentarios)"
$ws.Range("E9").Value = "planificada"
$ws.Range("F9").Value = "20 horas"
$ws.Range("G9").Value = 4
$ws.Range("H9").Value = "Alta"

$ws.Range("B10").Value = "H. usuario #6 y H. técnica #6"
$ws.Range("C10").Value = "Como un usuario activo, necesito recibir recompensas en tokens por`ninteractuar, con la finalidad de obtener un incentivo por mis actividades"
$ws.Range("D10").Value = "Sistema de tokenización por interacción"
$ws.Range("E10").Value = "planificada"
$ws.Range("F10").Value = "26 horas"
$ws.Range("G10").Value = 5
$ws.Range("H10").Value = "Alta"

$ws.Range("B11").Value = "H. usuario #8 y H. técnica #8"
$ws.Range("C11").Value = "Como un usuario, necesito enviar y recibir mensajes privados, con la finalidad de comunicarme de forma directa con otros usuarios"
$ws.Range("D11").Value = "Envío de mensajes privados"
$ws.Range("E11").Value = "planificada"
$ws.Range("F11").Value = "18 horas"
$ws.Range("G11").Value = 6
$ws.Range("H11").Value = "Alta"

$ws.Range("B12").Value = "H. usuario #13 y H. técnica #14"
$ws.Range("C12").Value = "Como un usuario, necesito verificar mi cuenta por email o teléfono, con la finalidad de mejorar la seguridad de mi cuenta y la confiabilidad dentro de la red social"
$ws.Range("D12").Value = "Verificación de cuenta"
$ws.Range("E12").Value = "planificada"
$ws.Range("F12").Value = "20 horas"
$ws.Range("G12").Value = 7
$ws.Range("H12").Value = "Alta"

$ws.Range("B13").Value = "H. usuario #7 y H. técnica #7"
$ws.Range("C13").Value = "Como un usuario, necesito buscar contenido en la red social, con la finalidad de descubrir publicaciones y usuarios que me interesen"
$ws.Range("D13").Value = "Búsqueda y exploración de contenidos"
$ws.Range("E13").Value = "planificada"
$ws.Range("F13").Value = "14 horas"
$ws.Range("G13").Value = 8
$ws.Range("H13").Value = "Media"

$ws.Range("B14").Value = "H. usuario #9 y H. técnica #10"
$ws.Range("C14").Value = "Como un usuario, necesito recibir notificaciones de actividades importantes  como reacciones a mis publicaciones o comentarios, con la finalidad de estar al tanto de las interacciones con mis publicaciones"
$ws.Range("D14").Value = "Notificaciones generales del sistema"
$ws.Range("E14").Value = "planificada"
$ws.Range("F14").Value = "20 horas"
$ws.Range("G14").Value = 9
$ws.Range("H14").Value = "Media"

$ws.Range("B15").Value = "H. usuario #12 y H. técnica #13"
$ws.Range("C15").Value = "Como un usuario, necesito ver una lista de mis amigos o seguidores, con la finalidad de gestionar mejor mis conexiones en la red social"
$ws.Range("D15").Value = "Lista de amigos y seguidores"
$ws.Range("E15").Value = "planificada"
$ws.Range("F15").Value = "12 horas"
$ws.Range("G15").Value = 10
$ws.Range("H15").Value = "Media"

$ws.Range("B16").Value = "H. usuario #11 y H. técnica #12"
$ws.Range("C16").Value = "Como un usuario, necesito  reportar contenido inapropiado como publicaciones ofensivas o spam, con la finalidad de mantener la comunidad segura y respetuosa"
$ws.Range("D16").Value = "Reporte de contenido inapropiado"
$ws.Range("E16").Value = "planificada"
$ws.Range("F16").Value = "16 horas"
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = "Media"

$ws.Range("B17").Value = "H. usuario #10 y H. técnica #11"
$ws.Range("C17").Value = "Como un usuario, necesito bloquear a otros usuarios, con la finalidad de evitar interacciones no deseadas o acosadoras"
$ws.Range("D17").Value = "Bloqueo de usuarios"
$ws.Range("E17").Value = "planificada"
$ws.Range("F17").Value = "18 horas"
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = "Baja"

$excel.ActiveWindow.Zoom = 85
$ws.Range("F17").Select()
